$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 198.33333
$ws.Range("I4").Value = 198.33333
$ws.Range("K4").Value = 198.33333
$ws.Range("M4").Value = -84.33332999999999

$ws.Range("H33").Value = 140.21428
$ws.Range("I33").Value = 112.69231
$ws.Range("K33").Value = 112.69231
$ws.Range("M33").Value = 116.30769

$ws.Range("H64").Value = 3994
$ws.Range("J64").Value = 3992.5
$ws.Range("L64").Value = 3992.5
$ws.Range("N64").Value = -4488.5

$ws.Range("H67").Value = 3994
$ws.Range("J67").Value = 3992.5
$ws.Range("L67").Value = 3992.5
$ws.Range("N67").Value = -5708.5

$ws.Range("H92").Value = 266.76923
$ws.Range("I92").Value = 218.63637
$ws.Range("K92").Value = 218.63637
$ws.Range("M92").Value = 1029.36363

$ws.Range("H112").Value = 1500.9584
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1531.4348
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4594.3044
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -6810.3044

$ws.Range("H121").Value = 1950.625
$ws.Range("J121").Value = 1950.625
$ws.Range("L121").Value = 5851.875
$ws.Range("N121").Value = -9345.875

$ws.Range("H124").Value = 180000
$ws.Range("J124").Value = 180000
$ws.Range("L124").Value = 180000
$ws.Range("N124").Value = -189820

$ws.Range("H137").Value = 2628.8857
$ws.Range("I137").Value = 1493.28
$ws.Range("J137").Value = 5467.9
$ws.Range("K137").Value = 4479.84
$ws.Range("L137").Value = 16403.7
$ws.Range("M137").Value = -1929.84
$ws.Range("N137").Value = -21503.7

$ws.Range("H138").Value = 2622.5688
$ws.Range("I138").Value = 1921.8572
$ws.Range("J138").Value = 2845.5227
$ws.Range("K138").Value = 5765.571599999999
$ws.Range("L138").Value = 8536.5681
$ws.Range("M138").Value = -625.5715999999993
$ws.Range("N138").Value = -18816.5681

$ws.Range("H141").Value = 5164.6665
$ws.Range("I141").Value = 5164.6665
$ws.Range("K141").Value = 15493.9995
$ws.Range("M141").Value = -10313.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8810.071
$ws.Range("I32").Value = 8810.071
$ws.Range("K32").Value = 8810.071
$ws.Range("M32").Value = -8523.071

$ws.Range("H61").Value = 2226.4443
$ws.Range("J61").Value = 2997.75
$ws.Range("L61").Value = 2997.75
$ws.Range("N61").Value = -3421.75

$ws.Range("H132").Value = 2906.9443
$ws.Range("I132").Value = 2318.9167
$ws.Range("J132").Value = 4083
$ws.Range("K132").Value = 6956.750100000001
$ws.Range("L132").Value = 12249
$ws.Range("M132").Value = -4426.750100000001
$ws.Range("N132").Value = -17309

$ws.Range("H136").Value = 2226.4443
$ws.Range("J136").Value = 2997.75
$ws.Range("L136").Value = 8993.25
$ws.Range("N136").Value = -14093.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 67254.39999999999
$ws.Range("J135").Value = 67254.39999999999
$ws.Range("L135").Value = 67254.39999999999
$ws.Range("N135").Value = -77394.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2167.8125
$ws.Range("I31").Value = 2197.9285
$ws.Range("K31").Value = 2197.9285
$ws.Range("M31").Value = -1902.9285

$ws.Range("H34").Value = 2167.8125
$ws.Range("I34").Value = 2197.9285
$ws.Range("K34").Value = 2197.9285
$ws.Range("M34").Value = -1995.9285

$ws.Range("H52").Value = 127125
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 127125
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 127125
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -127713

$ws.Range("H58").Value = 2165.9092
$ws.Range("I58").Value = 1714.3334
$ws.Range("K58").Value = 1714.3334
$ws.Range("M58").Value = -1511.3334

$ws.Range("H62").Value = 3273.25

$ws.Range("H65").Value = 3273.25

$ws.Range("H107").Value = 1087.6111
$ws.Range("I107").Value = 495.30768
$ws.Range("K107").Value = 495.30768
$ws.Range("M107").Value = 1424.69232

$ws.Range("H122").Value = 4080
$ws.Range("I122").Value = 3830
$ws.Range("J122").Value = 4580
$ws.Range("K122").Value = 11490
$ws.Range("L122").Value = 13740
$ws.Range("M122").Value = -9040
$ws.Range("N122").Value = -18640

$ws.Range("H136").Value = 2165.9092
$ws.Range("I136").Value = 1714.3334
$ws.Range("K136").Value = 5143.0002
$ws.Range("M136").Value = -2593.0002

$ws.Range("H141").Value = 59998.2
$ws.Range("J141").Value = 59998.2
$ws.Range("L141").Value = 59998.2
$ws.Range("N141").Value = -70358.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 222
$ws.Range("J33").Value = 222
$ws.Range("L33").Value = 1332
$ws.Range("N33").Value = -1898

$ws.Range("H56").Value = 18441.666
$ws.Range("I56").Value = 18441.666
$ws.Range("K56").Value = 18441.666
$ws.Range("M56").Value = -17911.666

$ws.Range("H131").Value = 1299.6
$ws.Range("J131").Value = 2498
$ws.Range("L131").Value = 7494
$ws.Range("N131").Value = -17574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.75
$ws.Range("I2").Value = 81.14286
$ws.Range("K2").Value = 81.14286
$ws.Range("M2").Value = 31.85714

$ws.Range("H126").Value = 1294.6364
$ws.Range("J126").Value = 1470.5714
$ws.Range("L126").Value = 4411.7142
$ws.Range("N126").Value = -9351.7142

$ws.Range("H132").Value = 2281
$ws.Range("J132").Value = 2989.4546
$ws.Range("L132").Value = 8968.363799999999
$ws.Range("N132").Value = -14028.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1350
$ws.Range("I46").Value = 1350
$ws.Range("K46").Value = 1350
$ws.Range("M46").Value = -1162

$ws.Range("H61").Value = 798.5
$ws.Range("I61").Value = 798.6667
$ws.Range("J61").Value = 798
$ws.Range("K61").Value = 798.6667
$ws.Range("L61").Value = 798
$ws.Range("M61").Value = -596.6667
$ws.Range("N61").Value = -1202

$ws.Range("H93").Value = 779.3333
$ws.Range("I93").Value = 998
$ws.Range("J93").Value = 670
$ws.Range("K93").Value = 998
$ws.Range("L93").Value = 670
$ws.Range("M93").Value = 250
$ws.Range("N93").Value = -3166

$ws.Range("H113").Value = 798.5
$ws.Range("I113").Value = 798.6667
$ws.Range("J113").Value = 798
$ws.Range("K113").Value = 798.6667
$ws.Range("L113").Value = 798
$ws.Range("M113").Value = 1371.3333
$ws.Range("N113").Value = -5138

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H107").Value = 562
$ws.Range("I107").Value = 562
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1686
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 234
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 815.94446
$ws.Range("I113").Value = 854.8
$ws.Range("J113").Value = 621.6667
$ws.Range("K113").Value = 2564.4
$ws.Range("L113").Value = 1865.0001
$ws.Range("M113").Value = -394.3999999999996
$ws.Range("N113").Value = -6205.0001

$ws.Range("H122").Value = 2062.6667
$ws.Range("I122").Value = 2070.5
$ws.Range("K122").Value = 6211.5
$ws.Range("M122").Value = -3761.5

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 2817.9
$ws.Range("I136").Value = 2860.5
$ws.Range("J136").Value = 2807.25
$ws.Range("K136").Value = 8581.5
$ws.Range("L136").Value = 8421.75
$ws.Range("M136").Value = -6031.5
$ws.Range("N136").Value = -13521.75
